# Automatic update of files.
# Bump the "Förändrad" date (column C) from 45293 to 45294 for data rows 2-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45293) {
        $cell.Value2 = 45294
    }
}
